{"js": "// Add the two intro paragraphs and a new \"10 The Next Room\" Heading1\n// paragraph in place of the single empty paragraph that used to sit\n// between the \"Write Up\" title and the trailing (still empty) Heading1\n// paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// paragraphs.items[0] = \"Write Up\" (Title)\n// paragraphs.items[1] = the empty paragraph we are replacing\n// paragraphs.items[2] = the trailing empty Heading1 paragraph (untouched)\nconst target = paragraphs.items[1];\n\n// Turn the empty paragraph into the first new body paragraph.\ntarget.insertText(\n  \"In this tutorial we will begin to start moving the hero toward the next room. We will be doing a bit more code in this section. Instead of placing a door into this next room, we will begin to have our code spawn the doors to exit the dungeon for us.\",\n  \"Replace\"\n);\n\n// Insert the second body paragraph right after it.\nconst second = target.insertParagraph(\n  \"So, if you would like to learn just a bit more about this then please join us for our brand-new article this week entitled:\",\n  \"After\"\n);\n\n// Insert the new Heading1 paragraph after the second paragraph.\nconst heading = second.insertParagraph(\"10 The Next Room\", \"After\");\nheading.styleBuiltIn = Word.BuiltInStyleName.heading1;\n\nawait context.sync();\n", "ps1": "# Replace the single empty paragraph that follows the \"Write Up\" title\n# with two body paragraphs and a new \"10 The Next Room\" Heading 1\n# paragraph. The existing (still empty) trailing Heading 1 paragraph is\n# left untouched.\n\n$d = $word.ActiveDocument\n\n# Paragraph 2 (1-based) is the empty paragraph right after the title.\n$p2 = $d.Paragraphs.Item(2)\n$r2 = $p2.Range\n$r2.Text = \"In this tutorial we will begin to start moving the hero toward the next room. We will be doing a bit more code in this section. Instead of placing a door into this next room, we will begin to have our code spawn the doors to exit the dungeon for us.\"\n$r2.InsertParagraphAfter()\n\n# Paragraph 3 is the new empty paragraph minted above.\n$p3 = $d.Paragraphs.Item(3)\n$r3 = $p3.Range\n$r3.Text = \"So, if you would like to learn just a bit more about this then please join us for our brand-new article this week entitled:\"\n$r3.InsertParagraphAfter()\n\n# Paragraph 4 is the next new empty paragraph; make it the Heading 1.\n$p4 = $d.Paragraphs.Item(4)\n$r4 = $p4.Range\n$r4.Text = \"10 The Next Room\"\n$p4.Style = \"Heading 1\"\n"}
